$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.745.20'
$ws.Range("E2").Value = '  +4.18%  '

$ws.Range("D3").Value = '3.071.95'
$ws.Range("E3").Value = '  +2.64%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.15'
$ws.Range("E5").Value = '  +2.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.44'
$ws.Range("E6").Value = '  +2.77%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.060.64'
$ws.Range("E8").Value = '  +2.70%  '

$ws.Range("E9").Value = '  +1.02%  '

$ws.Range("E10").Value = '  +5.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.64'
$ws.Range("E11").Value = '  +10.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.464'
$ws.Range("E12").Value = '  +2.06%  '

$ws.Range("E13").Value = '  +4.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.09'
$ws.Range("E14").Value = '  +4.41%  '

$ws.Range("E15").Value = '  +0.44%  '

$ws.Range("D16").Value = '3.578.15'
$ws.Range("E16").Value = '  +2.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.22'
$ws.Range("E17").Value = '  +0.39%  '

$ws.Range("D18").Value = '3.071.74'
$ws.Range("E18").Value = '  +2.69%  '

$ws.Range("D19").Value = '61.687.53'
$ws.Range("E19").Value = '  +4.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '445.80'
$ws.Range("E20").Value = '  +3.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.95'
$ws.Range("E21").Value = '  +2.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.731'
$ws.Range("E22").Value = '  +1.88%  '

$ws.Range("E23").Value = '  +4.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.69'
$ws.Range("E24").Value = '  +2.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.62'
$ws.Range("E25").Value = '  +0.82%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("E27").Value = '  +5.57%  '

$ws.Range("E28").Value = '  -0.05%  '

$ws.Range("E29").Value = '  +4.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.16'
$ws.Range("E30").Value = '  +5.34%  '

$ws.Range("E31").Value = '  +10.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.113'
$ws.Range("E32").Value = '  +14.43%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.74'
$ws.Range("E33").Value = '  +4.01%  '

$ws.Range("E34").Value = '  +4.13%  '

$ws.Range("D35").Value = '0.0₃0789'
$ws.Range("E35").Value = '  +3.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.02'
$ws.Range("E36").Value = '  +1.60%  '

$ws.Range("E37").Value = '  +4.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.99'
$ws.Range("E38").Value = '  +2.06%  '

$ws.Range("E39").Value = '  +9.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.79'
$ws.Range("E40").Value = '  +1.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '420.31'
$ws.Range("E41").Value = '  +4.74%  '

$ws.Range("D42").Value = '2.970.60'
$ws.Range("E42").Value = '  +7.81%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0368'
$ws.Range("E43").Value = '  +5.02%  '

$ws.Range("E44").Value = '  +9.91%  '

$ws.Range("E45").Value = '  +0.33%  '

$ws.Range("E46").Value = '  +5.69%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.55'
$ws.Range("E48").Value = '  +2.81%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.87'
$ws.Range("E49").Value = '  +0.40%  '

$ws.Range("E50").Value = '  +0.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.25'
$ws.Range("E51").Value = '  +3.64%  '
